$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 8340016
$ws.Range("I32").Value = 9500
$ws.Range("J32").Value = 11116855
$ws.Range("K32").Value = 9500
$ws.Range("L32").Value = 11116855
$ws.Range("M32").Value = -9174
$ws.Range("N32").Value = -11117507
$ws.Range("H40").Value = 3370.8215
$ws.Range("I40").Value = 2548
$ws.Range("J40").Value = 4642.4546
$ws.Range("K40").Value = 2548
$ws.Range("L40").Value = 4642.4546
$ws.Range("M40").Value = -2373
$ws.Range("N40").Value = -4992.4546
$ws.Range("H48").Value = 1257.1428
$ws.Range("J48").Value = 1257.1428
$ws.Range("L48").Value = 3771.4284
$ws.Range("N48").Value = -4355.428400000001
$ws.Range("H56").Value = 1257.1428
$ws.Range("J56").Value = 1257.1428
$ws.Range("L56").Value = 3771.4284
$ws.Range("N56").Value = -4839.428400000001
$ws.Range("H64").Value = 7684.467
$ws.Range("I64").Value = 5631.5
$ws.Range("K64").Value = 5631.5
$ws.Range("M64").Value = -5383.5
$ws.Range("H67").Value = 7684.467
$ws.Range("I67").Value = 5631.5
$ws.Range("K67").Value = 5631.5
$ws.Range("M67").Value = -4773.5
$ws.Range("H69").Value = 65178.6
$ws.Range("I69").Value = 8998.333000000001
$ws.Range("K69").Value = 26994.999
$ws.Range("M69").Value = -26120.999
$ws.Range("H72").Value = 65178.6
$ws.Range("I72").Value = 8998.333000000001
$ws.Range("K72").Value = 80984.997
$ws.Range("M72").Value = -76616.997
$ws.Range("H86").Value = 4052924.2
$ws.Range("I86").Value = 4240
$ws.Range("K86").Value = 4240
$ws.Range("M86").Value = -3117
$ws.Range("H89").Value = 4052924.2
$ws.Range("I89").Value = 4240
$ws.Range("K89").Value = 21200
$ws.Range("M89").Value = -15584
$ws.Range("H92").Value = 4808624.5
$ws.Range("I92").Value = 771.2105
$ws.Range("J92").Value = 17858512
$ws.Range("K92").Value = 771.2105
$ws.Range("L92").Value = 17858512
$ws.Range("M92").Value = 476.7895
$ws.Range("N92").Value = -17861008
$ws.Range("H97").Value = 2115.6667
$ws.Range("J97").Value = 2115.6667
$ws.Range("L97").Value = 6347.000100000001
$ws.Range("N97").Value = -7339.000100000001
$ws.Range("H98").Value = 2099.875
$ws.Range("I98").Value = 2239.8
$ws.Range("K98").Value = 2239.8
$ws.Range("M98").Value = -741.8000000000002
$ws.Range("H116").Value = 5000
$ws.Range("I116").Value = 5000
$ws.Range("K116").Value = 5000
$ws.Range("M116").Value = -1558
$ws.Range("H122").Value = 2099.875
$ws.Range("I122").Value = 2239.8
$ws.Range("K122").Value = 6719.400000000001
$ws.Range("M122").Value = -4269.400000000001
$ws.Range("H132").Value = 20620.309
$ws.Range("I132").Value = 1324
$ws.Range("K132").Value = 3972
$ws.Range("M132").Value = -1442
$ws.Range("H137").Value = 6648.8237
$ws.Range("I137").Value = 24725.666
$ws.Range("J137").Value = 2775.2144
$ws.Range("K137").Value = 74176.99800000001
$ws.Range("L137").Value = 8325.643199999999
$ws.Range("M137").Value = -71626.99800000001
$ws.Range("N137").Value = -13425.6432
$ws.Range("H141").Value = 3733.4
$ws.Range("I141").Value = 3733.4
$ws.Range("K141").Value = 11200.2
$ws.Range("M141").Value = -6020.200000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4104.92
$ws.Range("I61").Value = 3328.2666
$ws.Range("K61").Value = 3328.2666
$ws.Range("M61").Value = -3116.2666
$ws.Range("H102").Value = 25643264
$ws.Range("I102").Value = 2531.2
$ws.Range("J102").Value = 111112376
$ws.Range("K102").Value = 2531.2
$ws.Range("L102").Value = 111112376
$ws.Range("M102").Value = -909.1999999999998
$ws.Range("N102").Value = -111115620
$ws.Range("H134").Value = 64000
$ws.Range("J134").Value = 64000
$ws.Range("L134").Value = 64000
$ws.Range("N134").Value = -74140
$ws.Range("H136").Value = 4104.92
$ws.Range("I136").Value = 3328.2666
$ws.Range("K136").Value = 9984.799800000001
$ws.Range("M136").Value = -7434.799800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1773.5454
$ws.Range("I20").Value = 1876.875
$ws.Range("J20").Value = 1498
$ws.Range("K20").Value = 1876.875
$ws.Range("L20").Value = 1498
$ws.Range("M20").Value = -1629.875
$ws.Range("N20").Value = -1992
$ws.Range("H86").Value = 1666.7632
$ws.Range("I86").Value = 1106.091
$ws.Range("J86").Value = 2437.6875
$ws.Range("K86").Value = 1106.091
$ws.Range("L86").Value = 2437.6875
$ws.Range("M86").Value = 16.90900000000011
$ws.Range("N86").Value = -4683.6875
$ws.Range("H89").Value = 1666.7632
$ws.Range("I89").Value = 1106.091
$ws.Range("J89").Value = 2437.6875
$ws.Range("K89").Value = 5530.455
$ws.Range("L89").Value = 12188.4375
$ws.Range("M89").Value = 85.54500000000007
$ws.Range("N89").Value = -23420.4375
$ws.Range("H94").Value = 5052755.5
$ws.Range("I94").Value = 1473.3214
$ws.Range("J94").Value = 33339934
$ws.Range("K94").Value = 1473.3214
$ws.Range("L94").Value = 33339934
$ws.Range("M94").Value = -1022.3214
$ws.Range("N94").Value = -33340836
$ws.Range("H99").Value = 3388.5
$ws.Range("I99").Value = 3248
$ws.Range("J99").Value = 3810
$ws.Range("K99").Value = 3248
$ws.Range("L99").Value = 3810
$ws.Range("M99").Value = -1750
$ws.Range("N99").Value = -6806
$ws.Range("H105").Value = 3500.348
$ws.Range("I105").Value = 2923.7693
$ws.Range("J105").Value = 4249.9
$ws.Range("K105").Value = 2923.7693
$ws.Range("L105").Value = 4249.9
$ws.Range("M105").Value = -1176.7693
$ws.Range("N105").Value = -7743.9
$ws.Range("H107").Value = 1260.15
$ws.Range("I107").Value = 1255.421
$ws.Range("J107").Value = 1350
$ws.Range("K107").Value = 1255.421
$ws.Range("L107").Value = 1350
$ws.Range("M107").Value = 664.579
$ws.Range("N107").Value = -5190
$ws.Range("H134").Value = 3106.3877
$ws.Range("I134").Value = 2203.4614
$ws.Range("J134").Value = 6627.8
$ws.Range("K134").Value = 6610.3842
$ws.Range("L134").Value = 19883.4
$ws.Range("M134").Value = -4075.3842
$ws.Range("N134").Value = -24953.4
$ws.Range("H138").Value = 73496.766
$ws.Range("I138").Value = 70775
$ws.Range("J138").Value = 73723.586
$ws.Range("K138").Value = 70775
$ws.Range("L138").Value = 73723.586
$ws.Range("M138").Value = -65635
$ws.Range("N138").Value = -84003.586
$ws.Range("H140").Value = 80111.42999999999
$ws.Range("J140").Value = 80111.42999999999
$ws.Range("L140").Value = 80111.42999999999
$ws.Range("N140").Value = -90471.42999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 8169.143
$ws.Range("I16").Value = 9197.5
$ws.Range("J16").Value = 1999
$ws.Range("K16").Value = 9197.5
$ws.Range("L16").Value = 1999
$ws.Range("M16").Value = -8910.5
$ws.Range("N16").Value = -2573
$ws.Range("H31").Value = 2625.319
$ws.Range("I31").Value = 1724.5883
$ws.Range("K31").Value = 1724.5883
$ws.Range("M31").Value = -1429.5883
$ws.Range("H34").Value = 2625.319
$ws.Range("I34").Value = 1724.5883
$ws.Range("K34").Value = 1724.5883
$ws.Range("M34").Value = -1522.5883
$ws.Range("H43").Value = 34997.5
$ws.Range("J43").Value = 34997.5
$ws.Range("L43").Value = 34997.5
$ws.Range("N43").Value = -35365.5
$ws.Range("H58").Value = 5928.2
$ws.Range("I58").Value = 4082.875
$ws.Range("J58").Value = 7158.4165
$ws.Range("K58").Value = 4082.875
$ws.Range("L58").Value = 7158.4165
$ws.Range("M58").Value = -3879.875
$ws.Range("N58").Value = -7564.4165
$ws.Range("H74").Value = 55430.8
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 55430.8
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 55430.8
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -57178.8
$ws.Range("H77").Value = 55430.8
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 55430.8
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 166292.4
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -175028.4
$ws.Range("H95").Value = 10868.333
$ws.Range("J95").Value = 10868.333
$ws.Range("L95").Value = 10868.333
$ws.Range("N95").Value = -16360.333
$ws.Range("H99").Value = 10206386
$ws.Range("I99").Value = 2444308.8
$ws.Range("K99").Value = 2444308.8
$ws.Range("M99").Value = -2442810.8
$ws.Range("H101").Value = 34997.5
$ws.Range("J101").Value = 34997.5
$ws.Range("L101").Value = 34997.5
$ws.Range("N101").Value = -41487.5
$ws.Range("H105").Value = 11861
$ws.Range("I105").Value = 1296.3334
$ws.Range("K105").Value = 1296.3334
$ws.Range("M105").Value = 450.6666
$ws.Range("H107").Value = 8535.538
$ws.Range("I107").Value = 931.7273
$ws.Range("K107").Value = 931.7273
$ws.Range("M107").Value = 988.2727
$ws.Range("H113").Value = 8169.143
$ws.Range("I113").Value = 9197.5
$ws.Range("J113").Value = 1999
$ws.Range("K113").Value = 9197.5
$ws.Range("L113").Value = 1999
$ws.Range("M113").Value = -7027.5
$ws.Range("N113").Value = -6339
$ws.Range("H126").Value = 10206386
$ws.Range("I126").Value = 2444308.8
$ws.Range("K126").Value = 7332926.399999999
$ws.Range("M126").Value = -7330456.399999999
$ws.Range("H134").Value = 4901.523
$ws.Range("I134").Value = 3359.1304
$ws.Range("J134").Value = 6590.8096
$ws.Range("K134").Value = 10077.3912
$ws.Range("L134").Value = 19772.4288
$ws.Range("M134").Value = -7542.3912
$ws.Range("N134").Value = -24842.4288
$ws.Range("H136").Value = 5928.2
$ws.Range("I136").Value = 4082.875
$ws.Range("J136").Value = 7158.4165
$ws.Range("K136").Value = 12248.625
$ws.Range("L136").Value = 21475.2495
$ws.Range("M136").Value = -9698.625
$ws.Range("N136").Value = -26575.2495
$ws.Range("H139").Value = 59229.152
$ws.Range("J139").Value = 59229.152
$ws.Range("L139").Value = 59229.152
$ws.Range("N139").Value = -69509.152
$ws.Range("H141").Value = 28294.059
$ws.Range("J141").Value = 28294.059
$ws.Range("L141").Value = 28294.059
$ws.Range("N141").Value = -38654.059

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3221.88
$ws.Range("J5").Value = 4871.8
$ws.Range("L5").Value = 14615.4
$ws.Range("N5").Value = -14839.4
$ws.Range("H8").Value = 1055.3334
$ws.Range("I8").Value = 1055.3334
$ws.Range("K8").Value = 3166.0002
$ws.Range("M8").Value = -3027.0002
$ws.Range("H10").Value = 66
$ws.Range("J10").Value = 122
$ws.Range("L10").Value = 366
$ws.Range("N10").Value = -644
$ws.Range("H13").Value = 140.8
$ws.Range("J13").Value = 126.5
$ws.Range("L13").Value = 379.5
$ws.Range("N13").Value = -715.5
$ws.Range("H38").Value = 1092.9546
$ws.Range("J38").Value = 4479.8
$ws.Range("L38").Value = 13439.4
$ws.Range("N38").Value = -14133.4
$ws.Range("H40").Value = 48.235294
$ws.Range("I40").Value = 41.75
$ws.Range("J40").Value = 63.8
$ws.Range("K40").Value = 167
$ws.Range("L40").Value = 255.2
$ws.Range("M40").Value = -98
$ws.Range("N40").Value = -393.2
$ws.Range("H61").Value = 387
$ws.Range("I61").Value = 233.75
$ws.Range("J61").Value = 1000
$ws.Range("K61").Value = 701.25
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -486.25
$ws.Range("N61").Value = -3430
$ws.Range("H135").Value = 3221.88
$ws.Range("J135").Value = 4871.8
$ws.Range("L135").Value = 43846.2
$ws.Range("N135").Value = -48916.2
$ws.Range("H140").Value = 2944.0667
$ws.Range("J140").Value = 900
$ws.Range("L140").Value = 2700
$ws.Range("N140").Value = -13060

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 3846268.2
$ws.Range("I2").Value = 11.9
$ws.Range("J2").Value = 6250178.5
$ws.Range("K2").Value = 11.9
$ws.Range("L2").Value = 6250178.5
$ws.Range("M2").Value = 101.1
$ws.Range("N2").Value = -6250404.5
$ws.Range("H70").Value = 50040.668
$ws.Range("I70").Value = 77132.06
$ws.Range("K70").Value = 77132.06
$ws.Range("M70").Value = -76862.06
$ws.Range("H73").Value = 50040.668
$ws.Range("I73").Value = 77132.06
$ws.Range("K73").Value = 77132.06
$ws.Range("M73").Value = -76196.06
$ws.Range("H80").Value = 16725194
$ws.Range("I80").Value = 142043.88
$ws.Range("J80").Value = 27780628
$ws.Range("K80").Value = 142043.88
$ws.Range("L80").Value = 27780628
$ws.Range("M80").Value = -141045.88
$ws.Range("N80").Value = -27782624
$ws.Range("H83").Value = 16725194
$ws.Range("I83").Value = 142043.88
$ws.Range("J83").Value = 27780628
$ws.Range("K83").Value = 710219.4
$ws.Range("L83").Value = 138903140
$ws.Range("M83").Value = -705227.4
$ws.Range("N83").Value = -138913124

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 350
$ws.Range("J12").Value = 300
$ws.Range("L12").Value = 300
$ws.Range("N12").Value = -640
$ws.Range("H16").Value = 253.9
$ws.Range("I16").Value = 182.75
$ws.Range("J16").Value = 538.5
$ws.Range("K16").Value = 182.75
$ws.Range("L16").Value = 538.5
$ws.Range("M16").Value = -12.75
$ws.Range("N16").Value = -878.5
$ws.Range("H43").Value = 10000000
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("H46").Value = 852.6667
$ws.Range("I46").Value = 809.75
$ws.Range("J46").Value = 887
$ws.Range("K46").Value = 809.75
$ws.Range("L46").Value = 887
$ws.Range("M46").Value = -621.75
$ws.Range("N46").Value = -1263
$ws.Range("H55").Value = 1317.6316
$ws.Range("I55").Value = 1845.1666
$ws.Range("J55").Value = 413.2857
$ws.Range("K55").Value = 1845.1666
$ws.Range("L55").Value = 413.2857
$ws.Range("M55").Value = -1672.1666
$ws.Range("N55").Value = -759.2857
$ws.Range("H93").Value = 2354.8
$ws.Range("I93").Value = 2463.9
$ws.Range("K93").Value = 2463.9
$ws.Range("M93").Value = -1215.9
$ws.Range("H128").Value = 60416
$ws.Range("J128").Value = 60429
$ws.Range("L128").Value = 60429
$ws.Range("N128").Value = -70389
$ws.Range("H136").Value = 4759.3237
$ws.Range("I136").Value = 3016.7693
$ws.Range("K136").Value = 9050.3079
$ws.Range("M136").Value = -6500.3079

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 10501.667
$ws.Range("J4").Value = 10501.667
$ws.Range("L4").Value = 10501.667
$ws.Range("N4").Value = -10727.667
$ws.Range("H96").Value = 126409.5
$ws.Range("J96").Value = 1879
$ws.Range("L96").Value = 1879
$ws.Range("N96").Value = -4625
$ws.Range("H100").Value = 508.375
$ws.Range("I100").Value = 448.6154
$ws.Range("K100").Value = 897.2308
$ws.Range("M100").Value = -356.2308
$ws.Range("H107").Value = 1450.2
$ws.Range("I107").Value = 1083.6666
$ws.Range("K107").Value = 3250.9998
$ws.Range("M107").Value = -1330.9998
$ws.Range("H122").Value = 307422.16
$ws.Range("J122").Value = 5380.2144
$ws.Range("L122").Value = 16140.6432
$ws.Range("N122").Value = -21040.6432
$ws.Range("H124").Value = 398300
$ws.Range("J124").Value = 398300
$ws.Range("L124").Value = 398300
$ws.Range("N124").Value = -408120
$ws.Range("H126").Value = 1478.2222
$ws.Range("I126").Value = 1478.2222
$ws.Range("K126").Value = 4434.6666
$ws.Range("M126").Value = -1964.6666
$ws.Range("H135").Value = 68124.56
$ws.Range("J135").Value = 68124.56
$ws.Range("L135").Value = 68124.56
$ws.Range("N135").Value = -78264.56
$ws.Range("H137").Value = 66666
$ws.Range("J137").Value = 66666
$ws.Range("L137").Value = 66666
$ws.Range("N137").Value = -76866
